# Rename the existing sheet from "Sheet1" to "erosion", then add a new
# "soil" sheet after it and populate it with the soil measurement data,
# so the same workbook carries both the erosion dataset and the soil
# dataset (per commit message: "load erosion data and soil data in same
# excel file").

$wb = $excel.ActiveWorkbook
$erosion = $wb.ActiveSheet
$erosion.Name = "erosion"

# Add the new sheet right after "erosion" so tab order is erosion, soil.
$soil = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $erosion)
$soil.Name = "soil"

# Header row
$soil.Range("B1").Value = "PH值"
$soil.Range("C1").Value = "全氮含量(%)"
$soil.Range("D1").Value = "碳酸根离子"
$soil.Range("E1").Value = "硫酸根离子"
$soil.Range("F1").Value = "镁离子"
$soil.Range("G1").Value = "钾离子"
$soil.Range("H1").Value = "钠离子"
$soil.Range("I1").Value = "土壤电阻(Ω)"
$soil.Range("J1").Value = "站点"

# Data row
$soil.Range("A2").Value = 0
$soil.Range("B2").Value = 6.75
$soil.Range("C2").Value = 0.109
$soil.Range("D2").Value = 0.0126
$soil.Range("E2").Value = 0.0118
$soil.Range("F2").Value = 0.0017
$soil.Range("G2").Value = 0.0004
$soil.Range("H2").Value = 0.0049
$soil.Range("I2").Value = 32.9
$soil.Range("J2").Value = "沈阳站"

$soil.Range("A1:J2").Select()
